# Remove the "Quitar hora de las fechas" bullet item entirely.
# The bookmarkStart/bookmarkEnd ("_GoBack") that used to sit at the end of
# that paragraph must be preserved by merging it into the previous
# paragraph (the "Iconos en el menú... (no se alinean con el texto)" item).

$d = $word.ActiveDocument

# Locate the paragraph that contains the text to remove.
$target = $d.Content
$target.Find.Execute("Quitar hora de las fechas", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetPara = $target.Paragraphs.Item(1)

# Determine its paragraph index so we can grab the paragraph right before it.
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -eq $targetPara.Range.Start) {
        $targetIndex = $i
    }
}

$prevPara = $d.Paragraphs.Item($targetIndex - 1)

# Step 1: delete the paragraph mark that ends the previous paragraph. This
# merges the (currently empty-of-text) target paragraph into the previous
# one, while keeping the bookmarkStart/bookmarkEnd that live at the tail of
# the target paragraph (right before its own paragraph mark).
$mergeRange = $d.Range($prevPara.Range.End - 1, $prevPara.Range.End)
$mergeRange.Delete()

# Step 2: remove the now-merged "Quitar hora de las fechas" text itself,
# leaving the bookmark start/end (and the surviving paragraph mark) intact.
$target2 = $d.Content
$target2.Find.Execute("Quitar hora de las fechas", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target2.Delete()
